$wb = $excel.ActiveWorkbook

# --- AddOpportunity sheet -------------------------------------------------
$wsOpp = $wb.Worksheets.Item("AddOpportunity")

# AA2 keeps displaying "10" (its backing shared-string slot is simply
# renumbered upstream because of the sharedStrings reshuffle below).
$wsOpp.Range("AA2").Value = "10"

# --- AddContact sheet ------------------------------------------------------
$wsContact = $wb.Worksheets.Item("AddContact")

# I2 text is shortened - the "$10,000" retainer sentence is removed.
$wsContact.Range("I2").Value = "Opportunity Detail - SIC Code., Opportunity Detail - Tombstone Permission., Opportunity Description - Opportunity Description., Estimated Financials - Est. Transaction Size/Market Cap., Estimated Fees - Retainer, input zero if there's no Retainer fee., Referral Information - Referral Contact name is required., HL Internal Team - Team must include the following roles: Initiator, Seller, Principal, Manager, Associate(Optional), Analyst(Optional)., Legal Matters - Confidentiality Agreement, Conflicts Check - A Conflicts Check was completed more than 30 days ago. A new Conflicts Check must be completed., Administration - `"Women Led`" is required. Please update this field with the correct value, Administration - Date Engaged - Date of Executed Retainer or similar document., Approved FEIS form - Please complete and submit this form via the FEIS button., Opportunity Contacts - Add at least one Primary Opportunity Contact., Opportunity Contacts - Add at least one Billing Contact., Opportunity Contacts - Add at least one Contact with an approrpriate Role - confirm with FVA BUAs., Enter Transaction Type, Estimated Fees - Total Anticipated Revenue should be Greater Than or Equal to the Fee."

# J2 and K2 keep the same displayed text (only their backing shared-string slot changes upstream).
$wsContact.Range("J2").Value = "Error:, Opportunity Detail - Valuation Date."
$wsContact.Range("K2").Value = "Opportunity Detail - Client: Street Address., Opportunity Detail - Client: City Address., Opportunity Detail - Client: Postal Code., Opportunity Detail - Subject: Street Address., Opportunity Detail - Subject: City Address., Opportunity Detail - Subject: Postal Code Address., Opportunity Detail - Valuation Date."

# Row 2 shrinks now that I2 holds shorter wrapped text.
$wsContact.Rows.Item(2).RowHeight = 225

# --- View / selection state --------------------------------------------------
# Leave AddOpportunity selected first so the final Select() below (on
# AddContact) is what ends up as the active/tabSelected sheet, matching the
# workbook's original active tab.
$wsOpp.Range("AA4").Select()

$wsContact.Range("I2").Select()
